$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.690.79'
$ws.Range("E2").Value = '  -2.35%  '
$ws.Range("D3").Value = '2.367.51'
$ws.Range("E3").Value = '  -3.20%  '
$ws.Range("E4").Value = '  +0.07%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '571.06'
$c.Style = "Normal"
$ws.Range("E5").Value = '  -1.60%  '
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '139.59'
$c.Style = "Normal"
$ws.Range("E6").Value = '  -2.57%  '
$ws.Range("E7").Value = '  +0.01%  '
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '0.526'
$c.Style = "Normal"
$ws.Range("E8").Value = '  -0.92%  '
$ws.Range("D9").Value = '2.367.11'
$ws.Range("E10").Value = '  +1.10%  '
$ws.Range("E11").Value = '  +0.52%  '
$ws.Range("E12").Value = '  -2.14%  '
$ws.Range("E13").Value = '  -1.17%  '
$ws.Range("E14").Value = '  -1.99%  '
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '0.0000170'
$c.Style = "Normal"
$ws.Range("E15").Value = '  -1.56%  '
$ws.Range("D16").Value = '2.797.21'
$ws.Range("E16").Value = '  -0.12%  '
$ws.Range("D17").Value = '60.594.58'
$ws.Range("E17").Value = '  -2.34%  '
$ws.Range("D18").Value = '2.360.13'
$ws.Range("E18").Value = '  -2.56%  '
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '10.54'
$c.Style = "Normal"
$ws.Range("E19").Value = '  -2.69%  '
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '7.27'
$c.Style = "Normal"
$ws.Range("E20").Value = '  +1.42%  '
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '321.34'
$c.Style = "Normal"
$ws.Range("E21").Value = '  -2.16%  '
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '4.03'
$c.Style = "Normal"
$ws.Range("E22").Value = '  -1.49%  '
$ws.Range("E23").Value = '  +1.60%  '
$ws.Range("E24").Value = '  +0.07%  '
$ws.Range("E25").Value = '  -5.90%  '
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '64.63'
$c.Style = "Normal"
$ws.Range("E26").Value = '  -1.56%  '
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '8.59'
$c.Style = "Normal"
$ws.Range("E27").Value = '  -8.22%  '
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '571.27'
$c.Style = "Normal"
$ws.Range("E28").Value = '  -6.23%  '
$ws.Range("D29").Value = '2.504.23'
$ws.Range("D30").Value = '0.0₃0912'
$ws.Range("E30").Value = '  -3.79%  '
$ws.Range("E31").Value = '  -1.59%  '
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '1.35'
$c.Style = "Normal"
$ws.Range("E32").Value = '  -5.44%  '
$ws.Range("E33").Value = '  -2.43%  '
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '0.132'
$c.Style = "Normal"
$ws.Range("E34").Value = '  -5.75%  '
$ws.Range("E35").Value = '  -0.11%  '
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '4.62'
$c.Style = "Normal"
$ws.Range("E37").Value = '  -2.25%  '
$ws.Range("E38").Value = '  -3.15%  '
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '146.95'
$c.Style = "Normal"
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '18.13'
$c.Style = "Normal"
$ws.Range("E40").Value = '  -1.24%  '
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '5.08'
$c.Style = "Normal"
$ws.Range("E41").Value = '  -4.39%  '
$ws.Range("E42").Value = '  +0.03%  '
$ws.Range("E43").Value = '  -4.07%  '
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '40.85'
$c.Style = "Normal"
$ws.Range("E44").Value = '  -4.00%  '
$ws.Range("E45").Value = '  -4.69%  '
$ws.Range("D46").Value = '0.0₆0280'
$ws.Range("E46").Value = '  +19.22%  '
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '139.77'
$c.Style = "Normal"
$ws.Range("E47").Value = '  -2.23%  '
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '3.50'
$c.Style = "Normal"
$ws.Range("E48").Value = '  -3.39%  '
$ws.Range("E49").Value = '  -3.48%  '
$ws.Range("E50").Value = '  -3.97%  '
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '19.30'
$c.Style = "Normal"
$ws.Range("E51").Value = '  -0.74%  '
